$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing "Devreye al..." step text (row 17, column B) ---
$ws.Range("B17").Value = "Devreye al. Akım voltajı ölç ve yaz."

# --- Insert the new detailed calibration-change test steps (rows 18-23) ---
$ws.Range("B18").Value = "Akım voltajı değiştirme menüsüne gel."
$ws.Range("B19").Value = "Akım ve voltajı değiştir ve yaz."
$ws.Range("B20").Value = "Akım ve voltajı ölç ve yaz."
$ws.Range("B21").Value = "Ölçüm değeri ile ayar değerini karşılaştır."
$ws.Range("B22").Value = "Cihazı kapatıp aç. Ayar değeri duruyor mu kontrol et."
$ws.Range("B23").Value = "Akım ve voltajı ölç ve yaz."

# --- Move "Oto şarj test" further down, from A20 to A27 ---
$ws.Range("A20").ClearContents()
$ws.Range("A27").Value = "Oto şarj test"

# --- Widen column B to fit the longer descriptions ---
$ws.Columns.Item(2).ColumnWidth = 48

# --- Update the view: scroll so row 10 is at the top, and select B24 ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B24").Select()
